# Weekly update: insert a new record (week of 2023-10-24, serial 45223) at the
# top of the "Sandia" price history, pushing all existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 759..821 down to 760..822, creating a blank row 759.
$ws.Rows.Item(759).Insert()

# Populate the new row 759 with the latest weekly record.
$ws.Range("A759").Value = 9
$ws.Range("B759").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C759").Value = "Metropolitana"
$ws.Range("D759").Value = 45223
$ws.Range("D759").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E759").Value = 13
$ws.Range("F759").Value = 100112028
$ws.Range("G759").Value = "Sandia"
$ws.Range("H759").Value = "Sin especificar"
$ws.Range("I759").Value = "Primera"
$ws.Range("J759").Value = 430
$ws.Range("K759").Value = 700
$ws.Range("L759").Value = 800
$ws.Range("M759").Value = 750
$ws.Range("N759").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O759").Value = "Perú"
$ws.Range("P759").Value = 750
$ws.Range("Q759").Value = 1
$ws.Range("R759").Value = "Hortaliza"
